$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new K (column G) value
$kValues = @{
    2 = 3
    3 = 0
    5 = 0
    6 = 0
    7 = 2
    8 = 0
    9 = 0
    10 = 1
    11 = 0
    12 = 0
    13 = 0
    15 = 1
    16 = 1
    17 = 1
    18 = 0
    19 = 0
    20 = 1
    21 = 0
    22 = 0
    23 = 2
    24 = 0
    25 = 0
    26 = 1
    27 = 1
    28 = 1
    29 = 1
    30 = 2
    31 = 0
    32 = 1
    33 = 0
    34 = 1
    35 = 1
    36 = 0
    38 = 1
    39 = 0
    40 = 0
    41 = 0
    42 = 1
    43 = 0
    44 = 1
    45 = 0
    46 = 0
    47 = 0
    48 = 1
    49 = 2
    51 = 1
    52 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}

